$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells so numeric-looking strings (e.g. "1.00", "63.20")
# are preserved exactly instead of being normalized as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.444.00"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "1.589.39"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +1.10%  "
$ws.Range("D5").Value = "213.47"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("D8").Value = "24.14"
$ws.Range("E8").Value = "  +6.44%  "
$ws.Range("D9").Value = "0.251"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.817.24"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "1.589.05"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "28.445.16"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").Value = "63.20"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "230.29"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").Value = "7.48"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "9.34"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "151.66"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").Value = "1.14"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "3.17"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "1.401.38"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -8.07%  "
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +8.83%  "
$ws.Range("D40").Value = "0.541"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("D44").Value = "5.56"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").Value = "0.980"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "64.21"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "1.726.76"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "87.20"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0520"
$ws.Range("E51").Value = "  -1.00%  "
